$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 500.0
$ws.Range("C4").Value = 95.0
$ws.Range("C6").Value = 21956.0753832672
$ws.Range("C7").Value = 21576.0753832672
$ws.Range("C8").Value = 21297.393121769186
$ws.Range("C9").Value = 6840.0
$ws.Range("C11").Value = 3036.328800314118
$ws.Range("C13").Value = 18919.746582953085
$ws.Range("C14").Value = 18539.746582953085
$ws.Range("C15").Value = 12079.746582953088
$ws.Range("C16").Value = 11850.202937453087
$ws.Range("C17").Value = 11229.778937453088
$ws.Range("C21").Value = 215315.54665731726
$ws.Range("C22").Value = 211589.01965731726
$ws.Range("C23").Value = 208856.08025759773
$ws.Range("C24").Value = 67077.48599999998
$ws.Range("C27").Value = 185539.33282771683
$ws.Range("C28").Value = 181812.80582771683
$ws.Range("C29").Value = 118461.84682771686
$ws.Range("C30").Value = 116210.79263657428
$ws.Range("C31").Value = 110126.51161697431

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 2606.624999999999
$ws.Range("D6").Value = 29.3626816285427
$ws.Range("C7").Value = 2322.0
$ws.Range("D7").Value = -10.919292188174332
$ws.Range("C8").Value = 2830.0
$ws.Range("D8").Value = 8.569510382199242
$ws.Range("C9").Value = 2437.0
$ws.Range("D9").Value = -6.507456960629133
$ws.Range("D10").Value = -4.435812592912257
$ws.Range("D11").Value = 41.86927540401865
$ws.Range("C12").Value = 2858.333333333333
$ws.Range("D12").Value = 9.656484278840807

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 1737.7499999999995
$ws.Range("C7").Value = 2302.0
$ws.Range("D7").Value = 32.470148180118
$ws.Range("C8").Value = 1716.0
$ws.Range("D8").Value = -1.251618472162253
$ws.Range("C9").Value = 1934.0
$ws.Range("D9").Value = 11.29333908790105
$ws.Range("D10").Value = 32.29751114947493
$ws.Range("C11").Value = 2741.0
$ws.Range("D11").Value = 57.73270033088768
$ws.Range("C12").Value = 2257.0
$ws.Range("D12").Value = 29.880592720471906
$ws.Range("C13").Value = 1892.7142857142858
$ws.Range("D13").Value = 8.917524713813041

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 260.6624999999999
$ws.Range("D7").Value = -1.0214357646381471
$ws.Range("D8").Value = -42.83796096484916
$ws.Range("C9").Value = 122.0
$ws.Range("D9").Value = -53.196182803433544
$ws.Range("C10").Value = 176.33333333333331
$ws.Range("D10").Value = -32.351859844306944

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 260.6624999999999
$ws.Range("D7").Value = 26.600489138253533
$ws.Range("C8").Value = 226.0
$ws.Range("D8").Value = -13.297846832590004
$ws.Range("C9").Value = 278.0
$ws.Range("D9").Value = 6.651321152831752

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 260.6624999999999
$ws.Range("C3").Value = 243.99999999999994
$ws.Range("D3").Value = -68.79745520228902
$ws.Range("D9").Value = -60.10166402915646
$ws.Range("D10").Value = 305.12156524241135
$ws.Range("D11").Value = 47.31693281542229
$ws.Range("C12").Value = 122.0
$ws.Range("D16").Value = -60.10166402915646
$ws.Range("D17").Value = 305.12156524241135
$ws.Range("D18").Value = 47.31693281542229
$ws.Range("C19").Value = 122.0

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 695.0999999999998
$ws.Range("C5").Value = 734.0
$ws.Range("D5").Value = 5.596317076679663
$ws.Range("C6").Value = 878.0
$ws.Range("D6").Value = 26.312760753848426
$ws.Range("C7").Value = 998.0
$ws.Range("D7").Value = 43.57646381815573
$ws.Range("C8").Value = 890.0
$ws.Range("D8").Value = 28.039131060279157
$ws.Range("C9").Value = 875.0
$ws.Range("D9").Value = 25.88116817724071

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 695.0999999999998
$ws.Range("D5").Value = 234.3403826787514
$ws.Range("D6").Value = 234.37724375275826
